# Update "想去人数" (interest count) figures on the 展览, 演出 and 全部类型 sheets
# to reflect newly scraped totals (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 869
$wsExpo.Range("F3").Value = 13853
$wsExpo.Range("F4").Value = 13631
$wsExpo.Range("F11").Value = 59
$wsExpo.Range("F21").Value = 416
$wsExpo.Range("F25").Value = 99

# --- 演出 (Performances) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F6").Value = 168
$wsShow.Range("F7").Value = 1541

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 869
$wsAll.Range("F4").Value = 13853
$wsAll.Range("F5").Value = 13631
$wsAll.Range("F12").Value = 59
$wsAll.Range("F28").Value = 416
$wsAll.Range("F31").Value = 843
$wsAll.Range("F32").Value = 168
$wsAll.Range("F33").Value = 1541
$wsAll.Range("F37").Value = 99
